$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells keep their original text representation
# (values like "1.002", "0.05730", "29.998.57" are not valid numbers/
# would lose formatting if auto-converted), by forcing Text format before
# writing, then restoring the default "Normal" style so no stray number
# formatting is left behind.
$priceCells = $ws.Range("D2:D51")
$priceCells.NumberFormat = "@"

$ws.Range("D2").Value = "30.045.25"
$ws.Range("E2").Value = "  +3.78%  "
$ws.Range("D3").Value = "1.893.20"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "247.38"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.4972"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").Value = "44.78"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "0.2950"
$ws.Range("E9").Value = "  +5.93%  "
$ws.Range("D10").Value = "0.06638"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").Value = "1.893.53"
$ws.Range("E11").Value = "  +4.57%  "
$ws.Range("D12").Value = "16.99"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "0.07233"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "0.6787"
$ws.Range("E14").Value = "  +5.12%  "
$ws.Range("D15").Value = "86.09"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "4.848"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").Value = "30.015.83"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "0.000007924"
$ws.Range("E18").Value = "  +8.32%  "
$ws.Range("D19").Value = "0.9996"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "12.91"
$ws.Range("E20").Value = "  +5.23%  "
$ws.Range("D21").Value = "2.139.37"
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("D22").Value = "0.9982"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "4.766"
$ws.Range("E23").Value = "  +4.34%  "
$ws.Range("D24").Value = "5.677"
$ws.Range("E24").Value = "  +5.37%  "
$ws.Range("D25").Value = "9.215"
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").Value = "147.68"
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("D27").Value = "131.74"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").Value = "16.76"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("D29").Value = "1.961"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "1.361"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").Value = "4.222"
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").Value = "0.08761"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").Value = "3.941"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("D34").Value = "0.05091"
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("D35").Value = "1.119"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "0.7059"
$ws.Range("E36").Value = "  +4.59%  "
$ws.Range("D37").Value = "2.668"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").Value = "2.779"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("D39").Value = "2.232"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("D40").Value = "0.9476"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").Value = "0.01656"
$ws.Range("E41").Value = "  +3.89%  "
$ws.Range("D42").Value = "6.078"
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("D43").Value = "0.9978"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "103.41"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4212"
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("D46").Value = "7.481"
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("D47").Value = "0.1259"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").Value = "0.05730"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("D49").Value = "32.81"
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("D50").Value = "8.210"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "0.3732"
$ws.Range("E51").Value = "  +3.02%  "

$priceCells.Style = "Normal"
